$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate Species / Disease / Symptoms data (rows 2-25) to French
$ws.Range('A2').Value = 'volaille'
$ws.Range('B2').Value = 'La grippe aviaire'
$ws.Range('C2').Value = 'Fièvre, Mal de gorge, Toux, Maux de tête, Douleur musculaire, Respiration difficile, Conjonctivite (la_partie_intérieure_de_la_paupière)'
$ws.Range('A3').Value = 'bétail'
$ws.Range('B3').Value = 'fièvre de la vallée du Rift'
$ws.Range('C3').Value = 'Fièvre, léthargie, mort subite, écoulement nasal, salivation excessive, anorexie, diarrhée'
$ws.Range('A4').Value = 'bétail'
$ws.Range('B4').Value = 'Tuberculose bovine'
$ws.Range('C4').Value = 'état général de la maladie, perte de poids, mort subite'
$ws.Range('A5').Value = 'bétail'
$ws.Range('B5').Value = 'Peau grumeleuse'
$ws.Range('C5').Value = 'Nodules sur la peau, Peau lésée'
$ws.Range('A6').Value = 'bétail'
$ws.Range('B6').Value = 'Peste des Petits Ruminants'
$ws.Range('C6').Value = 'Fièvre, lésions buccales, diarrhée, mort subite'
$ws.Range('A7').Value = 'bétail'
$ws.Range('B7').Value = 'Diarrhée virale bovine'
$ws.Range('C7').Value = 'Fièvre, léthargie, anorexie, sécrétions oculaires, écoulement nasal, lésions buccales, diarrhée, réduction de la production'
$ws.Range('A8').Value = 'volaille'
$ws.Range('B8').Value = 'Newcastle'
$ws.Range('C8').Value = 'Sneezing, Nasal discharge, Cough, Diarrhea, Shudder, Drooping wings, Paralysis, Swelling of the tissues around the eyes and the neck, sudden death,  Reduction in production'
$ws.Range('A9').Value = 'volaille'
$ws.Range('B9').Value = 'Coryza infectieux'
$ws.Range('C9').Value = 'Tête ou visage enflé, éternuements, toux, sécrétions oculaires, écoulement nasal, anorexie, respiration difficile'
$ws.Range('A10').Value = 'volaille'
$ws.Range('B10').Value = 'Muguet'
$ws.Range('C10').Value = 'Anorexie, lésions buccales, croissance lente, yeux squameux, ailes tombantes, respiration difficile, respiration bruyante'
$ws.Range('A11').Value = 'bétail'
$ws.Range('B11').Value = 'Mastite (infection bactérienne)'
$ws.Range('C11').Value = 'Taille anormale, dureté de la mamelle, fièvre, changement de couleur du lait'
$ws.Range('A12').Value = 'bétail'
$ws.Range('B12').Value = 'Tuberculose'
$ws.Range('C12').Value = 'Abcès, toux, ganglions lymphatiques enflés, augmentation de la fréquence cardiaque'
$ws.Range('A13').Value = 'bétail'
$ws.Range('B13').Value = 'Gonfler'
$ws.Range('C13').Value = 'Flatulences, salivation excessive, gémissements, anorexie, vomissements'
$ws.Range('A14').Value = 'bétail'
$ws.Range('B14').Value = 'Encéphalomyélite équine'
$ws.Range('C14').Value = 'Fièvre, vision double, démarche irrégulière, frissons, lésions buccales, marche sans but'
$ws.Range('A15').Value = 'bétail'
$ws.Range('B15').Value = 'Septicémie hémorragique'
$ws.Range('C15').Value = 'Fièvre, respiration difficile, toux, sécrétions oculaires, écoulement nasal'
$ws.Range('A16').Value = 'bétail'
$ws.Range('B16').Value = 'Brucellose'
$ws.Range('C16').Value = 'Gonflement des testicules, bactéries localisées dans les articulations, Fièvre, Frissons, Anorexie, Transpiration, Léthargie'
$ws.Range('A17').Value = 'bétail'
$ws.Range('B17').Value = 'Clavelée'
$ws.Range('C17').Value = 'Fièvre, peau endommagée, inflammation de la muqueuse nasale, lésions cutanées, sécrétions oculaires, écoulement nasal, paupières enflées, léthargie, anorexie, paralysie'
$ws.Range('A18').Value = 'bétail'
$ws.Range('B18').Value = 'Pneumonie'
$ws.Range('C18').Value = 'Fièvre, anorexie, respiration difficile, sécrétions oculaires, écoulement nasal, salivation excessive, diarrhée'
$ws.Range('A19').Value = 'bétail'
$ws.Range('B19').Value = 'Fièvre aphteuse (FA)'
$ws.Range('C19').Value = 'Fièvre, lésions cutanées, peau endommagée, respiration difficile, salivation excessive, lésions buccales'
$ws.Range('A20').Value = 'bétail'
$ws.Range('B20').Value = 'Fièvre catarrhale du mouton'
$ws.Range('C20').Value = 'Fièvre, Lésions buccales, Respiration difficile, Langue violacée, Boiterie'
$ws.Range('A21').Value = 'bétail'
$ws.Range('B21').Value = 'Anaplasmose'
$ws.Range('C21').Value = 'Fièvre, pâleur autour des yeux, léthargie, perte de poids, réduction de la production, comportement agressif'
$ws.Range('A22').Value = 'bétail'
$ws.Range('B22').Value = 'Rage'
$ws.Range('C22').Value = 'Anorexie, prurit, boiterie, ténesme, salivation excessive, comportement agressif'
$ws.Range('A23').Value = 'volaille'
$ws.Range('B23').Value = 'Entérite nécrotique'
$ws.Range('C23').Value = 'Anorexie, Léthargie, Plumes pelucheuses, Yeux fermés, Diarrhée'
$ws.Range('A24').Value = 'volaille'
$ws.Range('B24').Value = 'Ascaris'
$ws.Range('C24').Value = 'Anorexie, Diarrhée, Croissance lente, Léthargie, Plumes ébouriffées, Perte de poids, changements de comportement'
$ws.Range('A25').Value = 'volaille'
$ws.Range('B25').Value = 'Variole aviaire'
$ws.Range('C25').Value = 'Paupières enflées, Yeux fermés, Lésions buccales, Perte de poids, Anorexie'

# Column widths to fit the longer French text
$ws.Columns.Item(1).ColumnWidth = 20.754
$ws.Columns.Item(2).ColumnWidth = 29.922
$ws.Columns.Item(3).ColumnWidth = 136.422

# Move the active selection to C25
$ws.Range("C25").Select()

Write-Output "Applied French translation edits"
